$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing In stock value for row 3 (T-shirts / Armany / Piece)
$ws.Cells.Item(3, 6).Value = 26

# Add new row 7: xvd
$ws.Cells.Item(7, 1).Value = "xvd"
$ws.Cells.Item(7, 2).Value = 1234
$ws.Cells.Item(7, 3).Value = "Armany"
$ws.Cells.Item(7, 4).Value = 12
$ws.Cells.Item(7, 5).Value = "Meter"
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = "22nd Dec, 2022"

# Add new row 8: Armany jeans
$ws.Cells.Item(8, 1).Value = "Armany jeans"
$ws.Cells.Item(8, 2).Value = 7869
$ws.Cells.Item(8, 3).Value = "Armany"
$ws.Cells.Item(8, 4).Value = 120
$ws.Cells.Item(8, 5).Value = "Piece"
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = "23rd Dec, 2022"
